# Scheduled runner update: refresh computed currentAveragePrice / Leve profit
# columns (H-N) across the per-job Sheets, mirroring the latest Universalis
# market-board pull. Values only; no structural/formula changes.
$wb = $excel.ActiveWorkbook

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2270956
$ws.Range("I132").Value = 440986.7
$ws.Range("K132").Value = 1322960.1
$ws.Range("M132").Value = -1320430.1

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 30410.75
$ws.Range("I2").Value = 32993.637
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 32993.637
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -32880.637
$ws.Range("N2").Value = -2225

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6487.9897
$ws.Range("I32").Value = 5159.011
$ws.Range("J32").Value = 19777.777
$ws.Range("K32").Value = 5159.011
$ws.Range("L32").Value = 19777.777
$ws.Range("M32").Value = -4872.011
$ws.Range("N32").Value = -20351.777

# ARM!row37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 7484.3
$ws.Range("J37").Value = 10015.2
$ws.Range("L37").Value = 10015.2
$ws.Range("N37").Value = -10561.2

# ARM!row41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 6343.5
$ws.Range("I41").Value = 3622.2
$ws.Range("J41").Value = 19950
$ws.Range("K41").Value = 3622.2
$ws.Range("L41").Value = 19950
$ws.Range("M41").Value = -3208.2
$ws.Range("N41").Value = -20778

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2798.8
$ws.Range("I45").Value = 1525.6364
$ws.Range("J45").Value = 6300
$ws.Range("K45").Value = 1525.6364
$ws.Range("L45").Value = 6300
$ws.Range("M45").Value = -1148.6364
$ws.Range("N45").Value = -7054

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 485.975
$ws.Range("I97").Value = 355.0357
$ws.Range("J97").Value = 791.5
$ws.Range("K97").Value = 355.0357
$ws.Range("L97").Value = 791.5
$ws.Range("M97").Value = 140.9643
$ws.Range("N97").Value = -1783.5

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2750
$ws.Range("I102").Value = 2125
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 2125
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -503
$ws.Range("N102").Value = -7244

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1188.75
$ws.Range("I110").Value = 421.72223
$ws.Range("J110").Value = 2569.4
$ws.Range("K110").Value = 421.72223
$ws.Range("L110").Value = 2569.4
$ws.Range("M110").Value = 1623.27777
$ws.Range("N110").Value = -6659.4

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 30410.75
$ws.Range("I116").Value = 32993.637
$ws.Range("J116").Value = 1999
$ws.Range("K116").Value = 32993.637
$ws.Range("L116").Value = 1999
$ws.Range("M116").Value = -30699.637
$ws.Range("N116").Value = -6587

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1432
$ws.Range("I122").Value = 1425.1904
$ws.Range("J122").Value = 1454
$ws.Range("K122").Value = 4275.5712
$ws.Range("L122").Value = 4362
$ws.Range("M122").Value = -1825.5712
$ws.Range("N122").Value = -9262

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 24418806
$ws.Range("I132").Value = 30055550
$ws.Range("J132").Value = 6996144
$ws.Range("K132").Value = 90166650
$ws.Range("L132").Value = 20988432
$ws.Range("M132").Value = -90164120
$ws.Range("N132").Value = -20993492

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 30410.75
$ws.Range("I3").Value = 32993.637
$ws.Range("J3").Value = 1999
$ws.Range("K3").Value = 32993.637
$ws.Range("L3").Value = 1999
$ws.Range("M3").Value = -32879.637
$ws.Range("N3").Value = -2227

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1983.1708
$ws.Range("I99").Value = 871.4286
$ws.Range("J99").Value = 2559.6296
$ws.Range("K99").Value = 871.4286
$ws.Range("L99").Value = 2559.6296
$ws.Range("M99").Value = 626.5714
$ws.Range("N99").Value = -5555.6296

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1975.75
$ws.Range("I105").Value = 1922
$ws.Range("J105").Value = 2014.1428
$ws.Range("K105").Value = 1922
$ws.Range("L105").Value = 2014.1428
$ws.Range("M105").Value = -175
$ws.Range("N105").Value = -5508.1428

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2500902.8
$ws.Range("I107").Value = 3334037
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 3334037
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -3332117
$ws.Range("N107").Value = -5340

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 918.94116
$ws.Range("I16").Value = 780.2222
$ws.Range("J16").Value = 1075
$ws.Range("K16").Value = 780.2222
$ws.Range("L16").Value = 1075
$ws.Range("M16").Value = -493.2222
$ws.Range("N16").Value = -1649

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1941681.8
$ws.Range("I31").Value = 3473701
$ws.Range("J31").Value = 6499.579
$ws.Range("K31").Value = 3473701
$ws.Range("L31").Value = 6499.579
$ws.Range("M31").Value = -3473406
$ws.Range("N31").Value = -7089.579

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1941681.8
$ws.Range("I34").Value = 3473701
$ws.Range("J34").Value = 6499.579
$ws.Range("K34").Value = 3473701
$ws.Range("L34").Value = 6499.579
$ws.Range("M34").Value = -3473499
$ws.Range("N34").Value = -6903.579

# CRP!row50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12390.8
$ws.Range("J50").Value = 12390.8
$ws.Range("L50").Value = 12390.8
$ws.Range("N50").Value = -13640.8

# CRP!row51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 33236.875
$ws.Range("J51").Value = 10982.5
$ws.Range("L51").Value = 10982.5
$ws.Range("N51").Value = -12454.5

# CRP!row55
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 7500
$ws.Range("I55").Value = 7500
$ws.Range("K55").Value = 7500
$ws.Range("M55").Value = -7185

# CRP!row59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16624.5
$ws.Range("J59").Value = 16624.5
$ws.Range("L59").Value = 16624.5
$ws.Range("N59").Value = -18914.5

# CRP!row60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18601.2
$ws.Range("J60").Value = 10376.5
$ws.Range("L60").Value = 10376.5
$ws.Range("N60").Value = -11398.5

# CRP!row61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 33236.875
$ws.Range("J61").Value = 10982.5
$ws.Range("L61").Value = 10982.5
$ws.Range("N61").Value = -11678.5

# CRP!row68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17910.889
$ws.Range("J68").Value = 17910.889
$ws.Range("L68").Value = 17910.889
$ws.Range("N68").Value = -19408.889

# CRP!row71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17910.889
$ws.Range("J71").Value = 17910.889
$ws.Range("L71").Value = 53732.667
$ws.Range("N71").Value = -61220.667

# CRP!row74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 17778.2
$ws.Range("J74").Value = 19610.777
$ws.Range("L74").Value = 19610.777
$ws.Range("N74").Value = -21358.777

# CRP!row77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 17778.2
$ws.Range("J77").Value = 19610.777
$ws.Range("L77").Value = 58832.33099999999
$ws.Range("N77").Value = -67568.33099999999

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 918.94116
$ws.Range("I113").Value = 780.2222
$ws.Range("J113").Value = 1075
$ws.Range("K113").Value = 780.2222
$ws.Range("L113").Value = 1075
$ws.Range("M113").Value = 1389.7778
$ws.Range("N113").Value = -5415

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 18809858
$ws.Range("I132").Value = 49527090
$ws.Range("J132").Value = 6995539
$ws.Range("K132").Value = 148581270
$ws.Range("L132").Value = 20986617
$ws.Range("M132").Value = -148578740
$ws.Range("N132").Value = -20991677

# LTW!row38
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2553137.5
$ws.Range("I136").Value = 2843041.5
$ws.Range("K136").Value = 8529124.5
$ws.Range("M136").Value = -8526574.5

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 15625911
$ws.Range("I126").Value = 20834098
$ws.Range("K126").Value = 62502294
$ws.Range("M126").Value = -62499824
